# Auto-applied numeric corrections to Leve profit calculation sheets
$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 1000
$ws.Range("I18").Value = 1000
$ws.Range("K18").Value = 1000
$ws.Range("M18").Value = -716
$ws.Range("H40").Value = 4000
$ws.Range("J40").Value = 6000
$ws.Range("L40").Value = 6000
$ws.Range("N40").Value = -6350
$ws.Range("H135").Value = 463.91666
$ws.Range("I135").Value = 322.68182
$ws.Range("K135").Value = 2904.13638
$ws.Range("M135").Value = -369.1363799999999
$ws.Range("H137").Value = 2611.7874
$ws.Range("J137").Value = 9500
$ws.Range("L137").Value = 28500
$ws.Range("N137").Value = -33600
$ws.Range("H141").Value = 35642.785
$ws.Range("I141").Value = 35642.785
$ws.Range("K141").Value = 106928.355
$ws.Range("M141").Value = -101748.355

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2640.5417
$ws.Range("I2").Value = 2522.9473
$ws.Range("K2").Value = 2522.9473
$ws.Range("M2").Value = -2409.9473
$ws.Range("H32").Value = 18522.555
$ws.Range("I32").Value = 3354
$ws.Range("K32").Value = 3354
$ws.Range("M32").Value = -3067
$ws.Range("H45").Value = 844393.4399999999
$ws.Range("I45").Value = 1685303
$ws.Range("K45").Value = 1685303
$ws.Range("M45").Value = -1684926
$ws.Range("H61").Value = 2691.804
$ws.Range("I61").Value = 2468.2222
$ws.Range("K61").Value = 2468.2222
$ws.Range("M61").Value = -2256.2222
$ws.Range("H74").Value = 1665.0526
$ws.Range("I74").Value = 1378.2858
$ws.Range("K74").Value = 1378.2858
$ws.Range("M74").Value = -504.2858000000001
$ws.Range("H77").Value = 1665.0526
$ws.Range("I77").Value = 1378.2858
$ws.Range("K77").Value = 6891.429
$ws.Range("M77").Value = -2523.429
$ws.Range("H87").Value = 0
$ws.Range("J87").Value = 0
$ws.Range("L87").Value = 0
$ws.Range("N87").ClearContents()
$ws.Range("H90").Value = 0
$ws.Range("J90").Value = 0
$ws.Range("L90").Value = 0
$ws.Range("N90").ClearContents()
$ws.Range("H107").Value = 45000
$ws.Range("J107").Value = 45000
$ws.Range("L107").Value = 45000
$ws.Range("N107").Value = -52680
$ws.Range("H116").Value = 2640.5417
$ws.Range("I116").Value = 2522.9473
$ws.Range("K116").Value = 2522.9473
$ws.Range("M116").Value = -228.9472999999998
$ws.Range("H132").Value = 3313.182
$ws.Range("I132").Value = 3313.182
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 9939.545999999998
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -7409.545999999998
$ws.Range("N132").ClearContents()
$ws.Range("H135").Value = 85238.5
$ws.Range("J135").Value = 85238.5
$ws.Range("L135").Value = 85238.5
$ws.Range("N135").Value = -95378.5
$ws.Range("H136").Value = 2691.804
$ws.Range("I136").Value = 2468.2222
$ws.Range("K136").Value = 7404.6666
$ws.Range("M136").Value = -4854.6666

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2640.5417
$ws.Range("I3").Value = 2522.9473
$ws.Range("K3").Value = 2522.9473
$ws.Range("M3").Value = -2408.9473
$ws.Range("H86").Value = 2692.1667
$ws.Range("I86").Value = 2280.6
$ws.Range("K86").Value = 2280.6
$ws.Range("M86").Value = -1157.6
$ws.Range("H89").Value = 2692.1667
$ws.Range("I89").Value = 2280.6
$ws.Range("K89").Value = 11403
$ws.Range("M89").Value = -5787
$ws.Range("H132").Value = 87654
$ws.Range("J132").Value = 87654
$ws.Range("L132").Value = 87654
$ws.Range("N132").Value = -97774
$ws.Range("H134").Value = 1439.3969
$ws.Range("I134").Value = 1238.695
$ws.Range("J134").Value = 4399.75
$ws.Range("K134").Value = 3716.085
$ws.Range("L134").Value = 13199.25
$ws.Range("M134").Value = -1181.085
$ws.Range("N134").Value = -18269.25

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 316.07693
$ws.Range("I22").Value = 316.07693
$ws.Range("K22").Value = 316.07693
$ws.Range("M22").Value = 33.92307
$ws.Range("H54").Value = 29999
$ws.Range("J54").Value = 29999
$ws.Range("L54").Value = 29999
$ws.Range("N54").Value = -31315
$ws.Range("H60").Value = 9665.416999999999
$ws.Range("J60").Value = 10199.2
$ws.Range("L60").Value = 10199.2
$ws.Range("N60").Value = -11221.2
$ws.Range("H97").Value = 29077.6
$ws.Range("J97").Value = 29077.6
$ws.Range("L97").Value = 29077.6
$ws.Range("N97").Value = -31059.6
$ws.Range("H105").Value = 1532.9412
$ws.Range("J105").Value = 1330
$ws.Range("L105").Value = 1330
$ws.Range("N105").Value = -4824
$ws.Range("H132").Value = 2401.3696
$ws.Range("I132").Value = 2330.9736
$ws.Range("K132").Value = 6992.9208
$ws.Range("M132").Value = -4462.9208
$ws.Range("H134").Value = 2480.162
$ws.Range("I134").Value = 2301.7812
$ws.Range("J134").Value = 3621.8
$ws.Range("K134").Value = 6905.3436
$ws.Range("L134").Value = 10865.4
$ws.Range("M134").Value = -4370.3436
$ws.Range("N134").Value = -15935.4
$ws.Range("H141").Value = 188697.75
$ws.Range("J141").Value = 188697.75
$ws.Range("L141").Value = 188697.75
$ws.Range("N141").Value = -199057.75

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 1766.1666
$ws.Range("J107").Value = 1899.25
$ws.Range("L107").Value = 5697.75
$ws.Range("N107").Value = -9537.75
$ws.Range("H122").Value = 3238.923
$ws.Range("J122").Value = 3323.6667
$ws.Range("L122").Value = 29913.0003
$ws.Range("N122").Value = -34813.0003
$ws.Range("H128").Value = 499997.5
$ws.Range("I128").Value = 499997.5
$ws.Range("K128").Value = 1499992.5
$ws.Range("M128").Value = -1495012.5
$ws.Range("H132").Value = 2107.1155
$ws.Range("I132").Value = 1349.3
$ws.Range("J132").Value = 2580.75
$ws.Range("K132").Value = 12143.7
$ws.Range("L132").Value = 23226.75
$ws.Range("M132").Value = -9613.699999999999
$ws.Range("N132").Value = -28286.75

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5476
$ws.Range("I70").Value = 5394.5
$ws.Range("J70").Value = 5639
$ws.Range("K70").Value = 5394.5
$ws.Range("L70").Value = 5639
$ws.Range("M70").Value = -5124.5
$ws.Range("N70").Value = -6179
$ws.Range("H73").Value = 5476
$ws.Range("I73").Value = 5394.5
$ws.Range("J73").Value = 5639
$ws.Range("K73").Value = 5394.5
$ws.Range("L73").Value = 5639
$ws.Range("M73").Value = -4458.5
$ws.Range("N73").Value = -7511
$ws.Range("H122").Value = 3831.8333
$ws.Range("I122").Value = 3831.8333
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 11495.4999
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -9045.499899999999
$ws.Range("N122").ClearContents()
$ws.Range("H132").Value = 2794.2766
$ws.Range("I132").Value = 2675.205
$ws.Range("K132").Value = 8025.615
$ws.Range("M132").Value = -5495.615

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 586.8570999999999
$ws.Range("J136").Value = 999
$ws.Range("L136").Value = 2997
$ws.Range("N136").Value = -8097
